$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.129.27'
$ws.Range("E2").Value = '  -1.24%  '
$ws.Range("D3").Value = '2.297.67'
$ws.Range("E3").Value = '  -2.11%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.28'
$ws.Range("E5").Value = '  -1.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.17'
$ws.Range("E6").Value = '  -0.57%  '
$ws.Range("E7").Value = '  -1.68%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.609'
$ws.Range("E9").Value = '  -1.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.00'
$ws.Range("E10").Value = '  -2.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0912'
$ws.Range("E11").Value = '  -1.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.42'
$ws.Range("E12").Value = '  +0.43%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.106'
$ws.Range("E13").Value = '  +0.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.973'
$ws.Range("E14").Value = '  -1.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.37'
$ws.Range("E15").Value = '  -3.73%  '
$ws.Range("D16").Value = '2.646.48'
$ws.Range("E16").Value = '  -2.03%  '
$ws.Range("D17").Value = '2.303.06'
$ws.Range("E17").Value = '  -3.54%  '
$ws.Range("D18").Value = '42.082.58'
$ws.Range("E18").Value = '  -1.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.61'
$ws.Range("E19").Value = '  -1.63%  '
$ws.Range("E20").Value = '  -0.62%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.87'
$ws.Range("E21").Value = '  -5.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.56'
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '259.61'
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.32'
$ws.Range("E24").Value = '  -0.32%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.89'
$ws.Range("E25").Value = '  +3.20%  '
$ws.Range("E26").Value = '  +0.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.98'
$ws.Range("E27").Value = '  -3.67%  '
$ws.Range("E28").Value = '  +2.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.74'
$ws.Range("E29").Value = '  -1.67%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.99'
$ws.Range("E30").Value = '  -1.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '164.91'
$ws.Range("E31").Value = '  -5.85%  '
$ws.Range("E32").Value = '  -0.20%  '
$ws.Range("E33").Value = '  -3.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.88'
$ws.Range("E34").Value = '  -3.04%  '
$ws.Range("E35").Value = '  -0.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.118'
$ws.Range("E36").Value = '  +5.55%  '
$ws.Range("E37").Value = '  +0.58%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0352'
$ws.Range("E38").Value = '  -1.47%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.90'
$ws.Range("E39").Value = '  +8.84%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.63'
$ws.Range("E40").Value = '  -3.81%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '99.57'
$ws.Range("E41").Value = '  +19.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.49'
$ws.Range("E42").Value = '  +1.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '70.73'
$ws.Range("E43").Value = '  +0.90%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.227'
$ws.Range("E44").Value = '  -1.94%  '
$ws.Range("E45").Value = '  +0.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.14'
$ws.Range("E46").Value = '  +1.97%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '113.86'
$ws.Range("E47").Value = '  -0.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '78.13'
$ws.Range("E48").Value = '  +7.38%  '
$ws.Range("E49").Value = '  -0.92%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.32'
$ws.Range("E50").Value = '  -3.39%  '
$ws.Range("E51").Value = '  +2.29%  '
